$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '55.799.94'
$ws.Range('E2').Value = '  +3.41%  '
$ws.Range('D3').Value = '2.448.65'
$ws.Range('E3').Value = '  +1.65%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '484.10'
$ws.Range('E5').Value = '  +4.05%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '143.67'
$ws.Range('E6').Value = '  +10.63%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.998'
$ws.Range('E7').Value = '  -0.15%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.505'
$ws.Range('E8').Value = '  +3.13%  '
$ws.Range('D9').Value = '2.451.02'
$ws.Range('E9').Value = '  +1.40%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '5.75'
$ws.Range('E10').Value = '  +8.80%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0963'
$ws.Range('E11').Value = '  +1.85%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.330'
$ws.Range('E12').Value = '  +5.23%  '
$ws.Range('E13').Value = '  +1.54%  '
$ws.Range('D14').Value = '2.874.18'
$ws.Range('E14').Value = '  +1.20%  '
$ws.Range('D15').Value = '55.914.37'
$ws.Range('E15').Value = '  +3.31%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '20.86'
$ws.Range('E16').Value = '  +6.74%  '
$ws.Range('E17').Value = '  +1.24%  '
$ws.Range('D18').Value = '2.451.09'
$ws.Range('E18').Value = '  +0.64%  '
$ws.Range('E19').Value = '  +6.67%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '10.01'
$ws.Range('E20').Value = '  +5.31%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '314.47'
$ws.Range('E21').Value = '  +1.42%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.998'
$ws.Range('E22').Value = '  +0.41%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.78'
$ws.Range('E23').Value = '  +7.52%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '58.10'
$ws.Range('E24').Value = '  +3.54%  '
$ws.Range('E25').Value = '  +5.89%  '
$ws.Range('E26').Value = '  -0.94%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.160'
$ws.Range('E27').Value = '  +3.55%  '
$ws.Range('D28').Value = '2.559.67'
$ws.Range('E28').Value = '  +0.47%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.66'
$ws.Range('E29').Value = '  +8.19%  '
$ws.Range('D30').Value = '0.0₃0772'
$ws.Range('E30').Value = '  +9.18%  '
$ws.Range('E31').Value = '  -0.10%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '147.50'
$ws.Range('E32').Value = '  +1.20%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '18.11'
$ws.Range('E33').Value = '  +2.31%  '
$ws.Range('E34').Value = '  +4.46%  '
$ws.Range('E35').Value = '  +2.33%  '
$ws.Range('E36').Value = '  +8.84%  '
$ws.Range('E37').Value = '  +3.95%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.845'
$ws.Range('E38').Value = '  +6.66%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '33.74'
$ws.Range('E39').Value = '  +3.46%  '
$ws.Range('E40').Value = '  +7.48%  '
$ws.Range('E41').Value = '  -0.10%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0549'
$ws.Range('E42').Value = '  +5.31%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.597'
$ws.Range('E43').Value = '  +0.58%  '
$ws.Range('E44').Value = '  +6.94%  '
$ws.Range('B45').Value = 'Stellar'
$ws.Range('C45').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0922'
$ws.Range('E45').Value = '  +4.72%  '
$ws.Range('B46').Value = 'Bittensor'
$ws.Range('C46').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '259.14'
$ws.Range('E46').Value = '  +12.78%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '10.18'
$ws.Range('E47').Value = '  +0.78%  '
$ws.Range('E48').Value = '  +12.24%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0227'
$ws.Range('E49').Value = '  +5.02%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '17.38'
$ws.Range('E50').Value = '  +5.24%  '
$ws.Range('D51').Value = '1.859.03'
$ws.Range('E51').Value = '  -3.10%  '
